$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9438162591565569, 0.7410130172864903, 0.9044778483062809, 3.045390320510702)"
$ws.Range("C2").Value = "JSU(-0.8452898259606749, 1.0909620827202438, 4.562415555171866, 4.745565372933688)"
$ws.Range("D2").Value = "NIG(0.6415280308414395, 0.3685070557919531, 2.4804543024132313, 2.8985106463598256)"
$ws.Range("E2").Value = "JSU(-1.0866085007192683, 1.1417333188554988, 3.2750257527091455, 5.575463199499277)"
